$wb = $excel.ActiveWorkbook

# --- Sheet "Layer0" ---
$ws1 = $wb.Worksheets.Item("Layer0")

$ws1.Range("B2").Value = -0.6910853737908073
$ws1.Range("C2").Value = 0.1252682626228348

$ws1.Range("B3").Value = 0.2891235804079033
$ws1.Range("C3").Value = -0.2907410849550568

$ws1.Range("B4").Value = 0.8260162381672392
$ws1.Range("C4").Value = -1.560745572121295

# --- Sheet "Layer1" ---
$ws2 = $wb.Worksheets.Item("Layer1")

$ws2.Range("B2").Value = -0.7379002882311573
$ws2.Range("C2").Value = -0.01531058539982979

$ws2.Range("B3").Value = 0.6079073878149142
$ws2.Range("C3").Value = -0.5090169693119415

$ws2.Range("B4").Value = -1.869709450317457
$ws2.Range("C4").Value = 0.7493145828825418
